{"js": "// This script rewrites the body of the document to match the target content.\n// The commit restructures the doc substantially: a reworded intro sentence,\n// merged overview/requirements bullets, a new Python implementation built\n// around validate_password()/re (replacing the old unittest-based code),\n// and a brand new \"Example Test Cases\" section with seven worked examples\n// (see commit message: \"Update Password Test Case SDET.docx from SDET sync\n// (code block from Test Lead)\").\n//\n// Rather than hand-chasing dozens of individual paragraph moves/renames,\n// we rebuild the body paragraph-by-paragraph from the final text, which is\n// the most reliable way to reproduce the diff exactly. Each array entry is\n// either a line of text or `null` for a blank paragraph (the document uses\n// blank paragraphs throughout as visual spacing between lines/sections).\nconst finalParagraphs = [\n  \"Here is the test code document for password validation including the overview, requirements, Python test code, and example test cases:\",\n  null,\n  \"# Test Cases for Password Validation \",\n  null,\n  \"## Overview\",\n  null,\n  \"This document covers test cases to validate password rules including minimum length of 8 characters, at least one numeric digit, and at least one special character.\",\n  null,\n  \"## Requirements\",\n  null,\n  \"- Password must be at least 8 characters long\",\n  \"- Password must contain at least one numeric digit \",\n  \"- Password must contain at least one special character\",\n  null,\n  \"## Python Test Code \",\n  null,\n  \"```python\",\n  \"import re\",\n  null,\n  \"def validate_password(password):\",\n  \"  if len(password) < 8:\",\n  \"    return False\",\n  \"  \",\n  \"  if not re.search(r\\\"\\\\d\\\", password):\",\n  \"    return False  \",\n  null,\n  \"  if not re.search(r\\\"[!@#$%^&*(),.?\\\\\\\":{}|<>]\\\", password):\",\n  \"    return False\",\n  null,\n  \"  return True\",\n  null,\n  \"print(validate_password(\\\"abc1$\\\")) # True\",\n  \"print(validate_password(\\\"abcd1234\\\")) # False\",\n  \"print(validate_password(\\\"abcd@xyz\\\")) # False\",\n  \"print(validate_password(\\\"abc1@def\\\")) # True\",\n  \"```\",\n  null,\n  \"## Example Test Cases\",\n  null,\n  \"**Test Case 1**\",\n  null,\n  \"Description: Password with less than 8 characters  \",\n  null,\n  \"Input: abc1$\",\n  null,\n  \"Expected Result: False\",\n  null,\n  \"**Test Case 2** \",\n  null,\n  \"Description: Password with 8 characters but no number\",\n  null,\n  \"Input: abcd@xyz\",\n  null,\n  \"Expected Result: False\",\n  null,\n  \"**Test Case 3**\",\n  null,\n  \"Description: Password with number but no special character\",\n  null,\n  \"Input: abcd1234\",\n  null,\n  \"Expected Result: False \",\n  null,\n  \"**Test Case 4**\",\n  null,\n  \"Description: Valid password with 8 characters, number, and special character\",\n  null,\n  \"Input: abc1@def\",\n  null,\n  \"Expected Result: True\",\n  null,\n  \"## Additional Test Cases\",\n  null,\n  \"**Test Case 5** \",\n  null,\n  \"Description: Valid password with more than 8 characters, number, and special character\",\n  null,\n  \"Input: MyPass123!\",\n  null,\n  \"Expected Result: True\",\n  null,\n  \"**Test Case 6**\",\n  null,\n  \"Description: Password with spaces\",\n  null,\n  \"Input: abcd 1234\",\n  null,\n  \"Expected Result: False\",\n  null,\n  \"**Test Case 7** \",\n  null,\n  \"Description: Password with allowed special characters \",\n  null,\n  \"Input: abcd!@#$\",\n  null,\n  \"Expected Result: True\"\n];\n\nconst body = context.document.body;\nbody.clear();\nawait context.sync();\n\nfor (const line of finalParagraphs) {\n  body.insertParagraph(line === null ? \"\" : line, Word.InsertLocation.end);\n}\nawait context.sync();\n\n// `clear()` leaves one leftover empty paragraph at the very start; remove it\n// now that the real content has been appended after it.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].delete();\nawait context.sync();\n", "ps1": "# This script rewrites the body of the document to match the target content.\n# The commit restructures the doc substantially: a reworded intro sentence,\n# merged overview/requirements bullets, a new Python implementation built\n# around validate_password()/re (replacing the old unittest-based code),\n# and a brand new \"Example Test Cases\" section with seven worked examples\n# (see commit message: \"Update Password Test Case SDET.docx from SDET sync\n# (code block from Test Lead)\").\n#\n# Rather than hand-chasing dozens of individual paragraph moves/renames,\n# we rebuild the body paragraph-by-paragraph from the final text, which is\n# the most reliable way to reproduce the diff exactly. Each array entry is\n# either a line of text or $null for a blank paragraph (the document uses\n# blank paragraphs throughout as visual spacing between lines/sections).\n$d = $word.ActiveDocument\n\n$finalParagraphs = @(\n    'Here is the test code document for password validation including the overview, requirements, Python test code, and example test cases:',\n    $null,\n    '# Test Cases for Password Validation ',\n    $null,\n    '## Overview',\n    $null,\n    'This document covers test cases to validate password rules including minimum length of 8 characters, at least one numeric digit, and at least one special character.',\n    $null,\n    '## Requirements',\n    $null,\n    '- Password must be at least 8 characters long',\n    '- Password must contain at least one numeric digit ',\n    '- Password must contain at least one special character',\n    $null,\n    '## Python Test Code ',\n    $null,\n    '```python',\n    'import re',\n    $null,\n    'def validate_password(password):',\n    '  if len(password) < 8:',\n    '    return False',\n    '  ',\n    '  if not re.search(r\"\\d\", password):',\n    '    return False  ',\n    $null,\n    '  if not re.search(r\"[!@#$%^&*(),.?\\\":{}|<>]\", password):',\n    '    return False',\n    $null,\n    '  return True',\n    $null,\n    'print(validate_password(\"abc1$\")) # True',\n    'print(validate_password(\"abcd1234\")) # False',\n    'print(validate_password(\"abcd@xyz\")) # False',\n    'print(validate_password(\"abc1@def\")) # True',\n    '```',\n    $null,\n    '## Example Test Cases',\n    $null,\n    '**Test Case 1**',\n    $null,\n    'Description: Password with less than 8 characters  ',\n    $null,\n    'Input: abc1$',\n    $null,\n    'Expected Result: False',\n    $null,\n    '**Test Case 2** ',\n    $null,\n    'Description: Password with 8 characters but no number',\n    $null,\n    'Input: abcd@xyz',\n    $null,\n    'Expected Result: False',\n    $null,\n    '**Test Case 3**',\n    $null,\n    'Description: Password with number but no special character',\n    $null,\n    'Input: abcd1234',\n    $null,\n    'Expected Result: False ',\n    $null,\n    '**Test Case 4**',\n    $null,\n    'Description: Valid password with 8 characters, number, and special character',\n    $null,\n    'Input: abc1@def',\n    $null,\n    'Expected Result: True',\n    $null,\n    '## Additional Test Cases',\n    $null,\n    '**Test Case 5** ',\n    $null,\n    'Description: Valid password with more than 8 characters, number, and special character',\n    $null,\n    'Input: MyPass123!',\n    $null,\n    'Expected Result: True',\n    $null,\n    '**Test Case 6**',\n    $null,\n    'Description: Password with spaces',\n    $null,\n    'Input: abcd 1234',\n    $null,\n    'Expected Result: False',\n    $null,\n    '**Test Case 7** ',\n    $null,\n    'Description: Password with allowed special characters ',\n    $null,\n    'Input: abcd!@#$',\n    $null,\n    'Expected Result: True'\n)\n\n# Clear the existing body down to a single empty paragraph. Deleting\n# Paragraphs.Item(1) repeatedly (i.e. always removing from the front) is the\n# reliable way to drain the document's content in this host.\nwhile ($d.Paragraphs.Count -gt 1) {\n    $d.Paragraphs.Item(1).Range.Delete()\n}\n$lead = $d.Paragraphs.Item(1).Range\n$lead.MoveEnd(1, -1) | Out-Null\n$lead.Text = \"\"\n\n# Append each target line as a new paragraph after the current last one.\n$cur = $d.Paragraphs.Item($d.Paragraphs.Count)\nforeach ($line in $finalParagraphs) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $d.Paragraphs.Item($d.Paragraphs.Count)\n    if ($line -ne $null) {\n        $r = $cur.Range\n        $r.MoveEnd(1, -1) | Out-Null\n        $r.Text = $line\n    }\n}\n\n# Remove the leftover leading empty paragraph left behind by the clear step.\n$d.Paragraphs.Item(1).Range.Delete()\n"}
